{"js": "// tdf#112118 fixture tweak: replace the leading page-break paragraph with\n// the \"_GoBack\" bookmark (start+end, empty), and remove that same bookmark\n// from the last paragraph (which becomes a plain empty paragraph).\n//\n// Before:\n//   P1: <w:r><w:br w:type=\"page\"/></w:r>\n//   P2: <w:pPr><w:sectPr>...</w:sectPr></w:pPr>   (unchanged)\n//   P3: <w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>\n// After:\n//   P1: <w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>\n//   P2: <w:pPr><w:sectPr>...</w:sectPr></w:pPr>   (unchanged)\n//   P3: <w:p/>  (now empty)\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\n\n// 1) Drop the existing \"_GoBack\" bookmark (lives in the last paragraph);\n//    this leaves that paragraph completely empty, matching the target.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Replace the first paragraph's content (the page-break run) with the\n//    \"_GoBack\" bookmark, using a raw OOXML swap so no stray empty run is\n//    left behind (unlike Range.clear() + Range.insertBookmark()).\nconst firstParagraphRange = firstParagraph.getRange();\nconst replacementOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p w:rsidR=\"001F344F\" w:rsidRDefault=\"001F344F\">\n<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n<w:bookmarkEnd w:id=\"0\"/>\n</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>`;\n\nfirstParagraphRange.insertOoxml(replacementOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# tdf#112118 fixture tweak: replace the leading page-break paragraph with\n# the \"_GoBack\" bookmark (start+end, empty), and remove that same bookmark\n# from the last paragraph (which becomes a plain empty paragraph).\n#\n# Before:\n#   P1: <w:r><w:br w:type=\"page\"/></w:r>\n#   P2: <w:pPr><w:sectPr>...</w:sectPr></w:pPr>   (unchanged)\n#   P3: <w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>\n# After:\n#   P1: <w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>\n#   P2: <w:pPr><w:sectPr>...</w:sectPr></w:pPr>   (unchanged)\n#   P3: <w:p/>  (now empty)\n\n$d = $word.ActiveDocument\n\n# 1) Drop the existing \"_GoBack\" bookmark (it lives in the last paragraph);\n#    removing it leaves that paragraph completely empty, matching the target.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) Remove the page-break character from the first paragraph (keep the\n#    paragraph mark itself) and drop the \"_GoBack\" bookmark at that now-empty\n#    spot instead.\n$firstParagraphRange = $d.Paragraphs(1).Range\n[void]$firstParagraphRange.MoveEnd(1, -1)\n$firstParagraphRange.Text = \"\"\n$d.Bookmarks.Add(\"_GoBack\", $firstParagraphRange)\n"}
